# Saldo_guide.xlsx update: refresh the extract date (2024-06-04 -> 2024-06-05)
# and correct a batch of balance values that had an extra leading
# thousands-digit typo (e.g. 11711.55 -> 711.55).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet to reflect the new export timestamp.
$ws.Name = "IClientBalance-20240605-095037-"

# 2) Bump every "G" column date (row 2..257) forward by one day
#    (serial 45447 -> 45448, i.e. 2024-06-04 -> 2024-06-05).
$lastRow = 257
for ($r = 2; $r -le $lastRow; $r++) {
    $cur = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value = $cur + 1
}

# 3) Correct the balances (column D) and their mirrored column H value
#    for the rows where the figure was overstated.
$corrections = @{
    5   = 711.55
    8   = 755.47
    15  = 1679.76
    17  = 253
    39  = 5177.53
    42  = 989.84
    57  = 1684.98
    59  = 450.27
    98  = 642.2
    103 = 792.95
    107 = 1054.59
    111 = 968.58
    131 = 449.98
    141 = 0.55
    143 = 0
    168 = 335.98
    226 = 179.5
    240 = 472.92
    245 = 4169.91
}

foreach ($row in $corrections.Keys) {
    $newValue = $corrections[$row]
    $ws.Cells.Item($row, 4).Value = $newValue
    $ws.Cells.Item($row, 8).Value = $newValue
}

Write-Output "Saldo_guide.xlsx updated: sheet renamed, dates advanced, $($corrections.Count) balances corrected."
